# Generate Report for Handback
# Re-runs the handback-status report generation with refreshed UUID-named
# files and updated timestamps:
#   c112864f-f6f7-44bb-abc2-93c28e1b4e25  ->  8231e5c8-dc3a-42dd-8f64-5e3650516065
#   e04279a5-149f-458b-8c7f-6d2d8006abd7  ->  ffff4da07da1-a5a2-408e-bc6e-4db1732007a6
# and consolidates the per-row handoff/handback xliff info onto the single
# refreshed xliff pair (both rows now reference the same generated files).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.md"
$ws1.Range("B2").Value = "e2e\8231e5c8-dc3a-42dd-8f64-5e3650516065.md"
$ws1.Range("G2").Value = "2016-08-18 00:58:50"

$ws1.Range("A3").Value = "ffff4da07da1-a5a2-408e-bc6e-4db1732007a6.md"
$ws1.Range("B3").Value = "e2e\ffff4da07da1-a5a2-408e-bc6e-4db1732007a6.md"
$ws1.Range("G3").Value = "2016-08-18 00:58:50"

foreach ($hl in $ws1.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\8231e5c8-dc3a-42dd-8f64-5e3650516065.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\ffff4da07da1-a5a2-408e-bc6e-4db1732007a6.md"
    }
}

# ---- zh-cn sheet ------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.md"
$ws2.Range("G2").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.9aca91942fa634ed5c2f517e2ba754c1693f820e.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-18 00:58:45"
$ws2.Range("I2").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.md"
$ws2.Range("J2").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.9aca91942fa634ed5c2f517e2ba754c1693f820e.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-18 00:59:05"

$ws2.Range("A3").Value = "ffff4da07da1-a5a2-408e-bc6e-4db1732007a6.md"
$ws2.Range("G3").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.9aca91942fa634ed5c2f517e2ba754c1693f820e.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-18 00:58:45"
$ws2.Range("I3").Value = "ffff4da07da1-a5a2-408e-bc6e-4db1732007a6.md"
$ws2.Range("J3").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.9aca91942fa634ed5c2f517e2ba754c1693f820e.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-18 00:59:05"

foreach ($hl in $ws2.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2' -or $addr -eq '$I$2') {
        $hl.TextToDisplay = "8231e5c8-dc3a-42dd-8f64-5e3650516065.md"
    } elseif ($addr -eq '$A$3' -or $addr -eq '$I$3') {
        $hl.TextToDisplay = "ffff4da07da1-a5a2-408e-bc6e-4db1732007a6.md"
    }
}

# ---- de-de sheet --------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.md"
$ws3.Range("G2").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.9aca91942fa634ed5c2f517e2ba754c1693f820e.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-18 00:58:50"
$ws3.Range("I2").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.md"
$ws3.Range("J2").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.9aca91942fa634ed5c2f517e2ba754c1693f820e.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-18 00:59:13"

$ws3.Range("A3").Value = "ffff4da07da1-a5a2-408e-bc6e-4db1732007a6.md"
$ws3.Range("G3").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.9aca91942fa634ed5c2f517e2ba754c1693f820e.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-18 00:58:50"
$ws3.Range("I3").Value = "ffff4da07da1-a5a2-408e-bc6e-4db1732007a6.md"
$ws3.Range("J3").Value = "8231e5c8-dc3a-42dd-8f64-5e3650516065.9aca91942fa634ed5c2f517e2ba754c1693f820e.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-18 00:59:13"

foreach ($hl in $ws3.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2' -or $addr -eq '$I$2') {
        $hl.TextToDisplay = "8231e5c8-dc3a-42dd-8f64-5e3650516065.md"
    } elseif ($addr -eq '$A$3' -or $addr -eq '$I$3') {
        $hl.TextToDisplay = "ffff4da07da1-a5a2-408e-bc6e-4db1732007a6.md"
    }
}
